# Daily attendance processing - 2025-11-23 14:48:01
# Reorders the "Recorded By" (column G) values so that any "System"/"system"
# tokens are moved to the end of the comma-separated list, while preserving
# the relative order of the remaining tokens (and of the System tokens
# themselves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text -split ", "
    $others = @()
    $systems = @()

    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systems += $p
        } else {
            $others += $p
        }
    }

    $newVal = @($others + $systems) -join ", "

    if ($newVal -ne $text) {
        $cell.Value2 = $newVal
    }
}
